# Refresh cryptos list prices / 1h volume deltas (scheduled GitHub Actions scrape).
# Column D prices that look like plain numbers are written with a leading
# apostrophe (Excel's text quote-prefix) so they stay text cells, matching the
# source data (which mixes thousands-dot-grouped prices with plain decimals).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.083.53"
$ws.Range("E2").Value = "  -0.43%  "

$ws.Range("D3").Value = "1.651.11"
$ws.Range("E3").Value = "  -0.57%  "

$ws.Range("D5").Value = "'218.16"
$ws.Range("E5").Value = "  -0.05%  "

$ws.Range("D6").Value = "'0.5283"
$ws.Range("E6").Value = "  +1.33%  "

$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").Value = "'0.2608"
$ws.Range("E8").Value = "  -2.13%  "

$ws.Range("D9").Value = "'0.06310"
$ws.Range("E9").Value = "  -0.19%  "

$ws.Range("D10").Value = "'20.33"
$ws.Range("E10").Value = "  -3.75%  "

$ws.Range("D11").Value = "'0.07739"
$ws.Range("E11").Value = "  +0.18%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.471"
$ws.Range("E12").Value = "  +0.94%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.656.86"
$ws.Range("E13").Value = "  -0.14%  "

$ws.Range("D14").Value = "'0.5452"
$ws.Range("E14").Value = "  -0.09%  "

$ws.Range("D15").Value = "0.0₅8104"
$ws.Range("E15").Value = "  -1.52%  "

$ws.Range("D16").Value = "'65.04"
$ws.Range("E16").Value = "  +0.13%  "

$ws.Range("D17").Value = "26.095.33"
$ws.Range("E17").Value = "  -0.57%  "

$ws.Range("E18").Value = "  -0.24%  "

$ws.Range("D19").Value = "'4.549"
$ws.Range("E19").Value = "  -2.43%  "

$ws.Range("D20").Value = "'193.33"
$ws.Range("E20").Value = "  +0.08%  "

$ws.Range("E21").Value = "  -1.16%  "

$ws.Range("D22").Value = "'5.985"
$ws.Range("E22").Value = "  -1.84%  "

$ws.Range("E23").Value = "  -0.29%  "

$ws.Range("D24").Value = "'140.07"
$ws.Range("E24").Value = "  +1.12%  "

$ws.Range("D25").Value = "'0.1240"
$ws.Range("E25").Value = "  -0.19%  "

$ws.Range("D26").Value = "'7.234"
$ws.Range("E26").Value = "  +0.30%  "

$ws.Range("D27").Value = "'16.15"
$ws.Range("E27").Value = "  -0.16%  "

$ws.Range("D28").Value = "'1.433"
$ws.Range("E28").Value = "  +0.93%  "

$ws.Range("D29").Value = "'0.05906"
$ws.Range("E29").Value = "  -1.66%  "

$ws.Range("D30").Value = "'1.279"
$ws.Range("E30").Value = "  -0.36%  "

$ws.Range("D31").Value = "'3.498"
$ws.Range("E31").Value = "  -2.13%  "

$ws.Range("D32").Value = "'3.237"
$ws.Range("E32").Value = "  -2.75%  "

$ws.Range("E33").Value = "  -6.42%  "

$ws.Range("D34").Value = "'2.411"
$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("D35").Value = "'0.9409"
$ws.Range("E35").Value = "  -4.14%  "

$ws.Range("D36").Value = "'2.751"
$ws.Range("E36").Value = "  -1.08%  "

$ws.Range("D37").Value = "'0.5670"
$ws.Range("E37").Value = "  -4.00%  "

$ws.Range("D38").Value = "'0.01602"
$ws.Range("E38").Value = "  +1.06%  "

$ws.Range("D39").Value = "'5.842"
$ws.Range("E39").Value = "  -1.72%  "

$ws.Range("D40").Value = "'0.8432"
$ws.Range("E40").Value = "  -2.44%  "

$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("D42").Value = "1.007.05"
$ws.Range("E42").Value = "  -3.19%  "

$ws.Range("D43").Value = "'100.68"
$ws.Range("E43").Value = "  +1.07%  "

$ws.Range("D44").Value = "1.799.93"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D45").Value = "'56.84"
$ws.Range("E45").Value = "  -0.41%  "

$ws.Range("E46").Value = "  -1.85%  "

$ws.Range("E47").Value = "  +0.06%  "

$ws.Range("E48").Value = "  +1.59%  "

$ws.Range("D49").Value = "'1.482"
$ws.Range("E49").Value = "  +1.78%  "

$ws.Range("D50").Value = "'0.05150"
$ws.Range("E50").Value = "  -0.59%  "

$ws.Range("D51").Value = "'7.826"
$ws.Range("E51").Value = "  -3.27%  "
